$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume snapshot (and fix the
# PaxDollar/Quant row ordering) per the latest GitHub Actions run.
# Values that look like plain numbers (e.g. "313.31", "1.012") are
# entered with a leading apostrophe so Excel keeps them as literal text
# instead of re-parsing them as numeric values.

$ws.Range("D2").Value = "27.855.87"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "1.874.09"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'313.31"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").Value = "'0.4829"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'0.3816"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("D9").Value = "'0.07367"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'0.9418"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "'21.03"
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("D12").Value = "'0.07786"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.894.70"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'5.518"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "'6.619"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'0.000008855"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "27.878.82"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "'14.87"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'5.125"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "2.121.37"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'10.87"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'1.948"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'157.43"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").Value = "'18.57"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'2.046"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("D29").Value = "'115.98"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'4.976"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'0.08886"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "'0.7697"
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("D35").Value = "'4.655"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").Value = "'2.741"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "'0.02043"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'0.5626"
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").Value = "'0.05371"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "'3.001"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'7.049"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "'8.563"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "'0.1532"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "'0.4881"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").Value = "'10.70"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.012"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'104.90"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "'0.06124"
$ws.Range("E51").Value = "  +0.83%  "
